# LocalAirData.xlsx update
#
# 1. Sheet1!E2 ("ToLocation" data) changes from "bost" to "miami".
# 2. The active sheet's view/selection moves from D11 (scrolled so B1 is the
#    top-left cell) to D9, with the sheet scrolled back to the default
#    top-left (A1), so the topLeftCell override is no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the destination city value.
$ws.Range("E2").Value = "miami"

# Make sure the sheet is scrolled to the top-left (A1) before selecting,
# so no topLeftCell override is persisted, then select D9.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D9").Select()
